$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: fix the test date (was 12/30/20, now 12/29/20) ---
$ws.Range("C12").Value = 44194

# --- Row 13: fix the test date (was 12/31/20, now 12/29/20) ---
$ws.Range("C13").Value = 44194

# --- Row 14: fill in the new "Egész törlés" (clear-all button) test case,
#     mirroring the row above it (row 13, the "Törlés gomb" / delete button test) ---
$ws.Range("B14").Value = "Szicsák Bence"
$ws.Range("C14").Value = 44194
$ws.Range("D14").Value = "x"
$ws.Range("E14").Value = "Igen"
$ws.Range("F14").Value = "Egész törles"
$ws.Range("G14").Value = "Sikerült"

# Row 13's "F" cell (the button-name column) uses a smaller font for the long
# label - match that same look for the new row's button-name cell.
$dstF = $ws.Range("F14")
$dstF.Font.Name = "Times New Roman"
$dstF.Font.Size = 8

# Row 11's "G" cell (Sikeresség/result column) is highlighted green for this
# tester - carry the same highlight onto the new row's result cell.
$dstG = $ws.Range("G14")
$dstG.Font.Name = "Times New Roman"
$dstG.Font.Size = 12
$dstG.Font.Color = 5287936

# --- Selection ends on J14, matching where the editor left the cursor ---
$ws.Range("J14").Select()
